$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing B:E block one column right (B->C, C->D, D->E, E->F).
# Excel's native column insert preserves the original custom width
# (14.83203125) on the surviving columns and auto-adjusts formula
# references, so it reproduces the XML shift far more faithfully than
# re-writing every cell by hand.
$ws.Columns("B").Insert() | Out-Null

# --- Row 2: drop the now-stray old E2 value (28 moved to F2, no longer needed) and set the new measurement ---
$ws.Range("F2").ClearContents()
$ws.Range("D2").Value = 3589

# --- Row 4: new label in A4 ---
$ws.Range("A4").Value = "T="

# --- Row 5: A5 becomes the numeric temperature, B5 becomes the "CB" label, formulas reference $A5 instead of a fixed cell ---
$ws.Range("A5").Value = 28
$ws.Range("B5").Value = "CB"
$ws.Range("D5").Formula = '=ABS(($D$2/1000*23.156-12.736)-(1+0.018*($A5-24)))'
$ws.Range("E5").Formula = '=ABS((-3.213*($D$2/1000)-4.093)/(1-0.009733*($D$2/1000)-0.01205*($A5)))'
$ws.Range("F5").Formula = '=ABS(2.246-5.239*($D$2/1000)*(1+0.018*($A5-24))-0.06756*($D$2/1000)*($D$2/1000)*((1+0.018*($A5-24))*(1+0.018*($A5-24))))'

# --- Row 6: new row, same shape, temperature 20 ---
$ws.Range("A6").Value = 20
$ws.Range("B6").Value = "CB"
$ws.Range("C6").Formula = "=0"
$ws.Range("D6").Formula = '=ABS(($D$2/1000*23.156-12.736)-(1+0.018*($A6-24)))'
$ws.Range("E6").Formula = '=ABS((-3.213*($D$2/1000)-4.093)/(1-0.009733*($D$2/1000)-0.01205*($A6)))'
$ws.Range("F6").Formula = '=ABS(2.246-5.239*($D$2/1000)*(1+0.018*($A6-24))-0.06756*($D$2/1000)*($D$2/1000)*((1+0.018*($A6-24))*(1+0.018*($A6-24))))'
$ws.Range("C6:F6").NumberFormat = "0.000"

# --- Selection matches the authored state ---
$ws.Range("E6").Select() | Out-Null
